$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValue = 93.61576470588236

for ($r = 35; $r -le 63; $r++) {
    $ws.Cells.Item($r, 9).Value = $newValue
}
